# Fruta / hortaliza, semanal
# Update weekly price records (rows 2-7 and 10-14) with refreshed
# date / quality / volume / price values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = 44253; I = "Segunda"; J = 1000; K = 800; L = 900; M = 850; P = 850 },
    @{ Row = 3;  D = 44253; I = "Tercera"; J = 800;  K = 600; L = 700; M = 650; P = 650 },
    @{ Row = 4;  D = 44267; I = "Tercera"; J = 400;  K = 500; L = 600; M = 550; P = 550 },
    @{ Row = 5;  D = 44174; I = "Segunda"; J = 800;  K = 450; L = 500; M = 475; P = 475 },
    @{ Row = 6;  D = 44174; I = "Tercera"; J = 1200; K = 250; L = 350; M = 300; P = 300 },
    @{ Row = 7;  D = 44224; I = "Segunda"; J = 800;  K = 850; L = 900; M = 875; P = 875 },
    @{ Row = 10; D = 44210; I = "Segunda"; J = 900;  K = 600; L = 700; M = 650; P = 650 },
    @{ Row = 11; D = 44278; I = "Segunda"; J = 700;  K = 600; L = 700; M = 650; P = 650 },
    @{ Row = 12; D = 44278; I = "Tercera"; J = 400;  K = 500; L = 600; M = 550; P = 550 },
    @{ Row = 13; D = 44474; I = "Segunda"; J = 200;  K = 600; L = 700; M = 650; P = 650 },
    @{ Row = 14; D = 44229; I = "Segunda"; J = 760;  K = 550; L = 600; M = 575; P = 575 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("D$r").Value = $u.D
    $ws.Range("I$r").Value = $u.I
    $ws.Range("J$r").Value = $u.J
    $ws.Range("K$r").Value = $u.K
    $ws.Range("L$r").Value = $u.L
    $ws.Range("M$r").Value = $u.M
    $ws.Range("P$r").Value = $u.P
}
